$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: clone the format of an already-styled reference cell (so the exact
# named cell style + border combination is reused, matching what Excel does
# when a user fills/drags a styled "Gantt" cell across the row) and then set
# the cell's text.
function Set-StyledCell($addr, $srcAddr, $text) {
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($addr)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $dst.Value = $text
}

$excel.CutCopyMode = $false

# --- Row 13 ("PYTHON - Bases de datos - Meter datos"): M13, N13, O13 move
#     from Retrasado/Pendiente to Realizado ---
Set-StyledCell "M13" "K13" "Realizado"
Set-StyledCell "N13" "K13" "Realizado"
Set-StyledCell "O13" "K13" "Realizado"

# --- Row 14 ("HTML"): O14 becomes Realizado, which means the following
#     "Puesta en comun" cell (P14) now caps a green bar instead of an
#     orange one, and the bar overruns into Q14/R14 (Retrasado) ---
Set-StyledCell "O14" "K13" "Realizado"
Set-StyledCell "P14" "E7"  "Puesta en común"
Set-StyledCell "Q14" "C4"  "Retrasado"
Set-StyledCell "R14" "C4"  "Retrasado"

# --- Row 15 ("CSS"): work starts (commit: "inicio de css en pantalla
#     clientes") - P15 becomes Realizado, Q15/R15 become Retrasado ---
Set-StyledCell "P15" "K13" "Realizado"
Set-StyledCell "Q15" "C4"  "Retrasado"
Set-StyledCell "R15" "C4"  "Retrasado"

$excel.CutCopyMode = $false

# --- Move the "today" line connector further right/down ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 1429.4765625
$shp.Top = 80.25
$shp.Height = 206.25

# --- Selection moves to P35 ---
$ws.Range("P35").Select() | Out-Null
